# OPAR Ordnance.xlsx update:
#  - "CURRENT AS OF:" value advances from SPINS version D2.1 to D2.2
#  - GBU-31(V) 3/B (row 9) expends 2 rounds under the new D2.2 column (G)
#  - GBU-38 (row 12) expends 4 rounds under the new D2.2 column (G)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Current-as-of SPINS version bump: D2.1 -> D2.2
$ws.Range("R1").Value = "D2.2"

# Newly reported ordnance usage for SPINS version D2.2 (column G)
$ws.Range("G9").Value = 2
$ws.Range("G12").Value = 4

# Reflect the updated selection left behind by the author
$ws.Range("R1:R2").Select()
